# "Generate Report for handoff"
#
# Updates the localization-status workbook so the zh-cn and de-de sheets
# reflect that a handoff package has been generated and is ready:
#   - Status moves from "Handoff transform failed" to "Ready for handoff"
#   - A "Latest Handoff File" hyperlink cell is populated with the newly
#     generated xlf handoff file for each locale
#   - "Latest Handoff Datetime" is stamped with the generation time
#   - "Handoff Reason" moves from "Ignored" to "Include"

$wb = $excel.ActiveWorkbook

$repoBase = "https://github.com/OpenLocalizationTest/oltest/blob/b76cd49c4b4099e031447afa5290f0430fe9375c/e2e"

# The Overview sheet mirrors each locale's Status in columns B/C, so it
# also needs to read "Ready for handoff" once the handoff succeeds.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"

function Update-LocaleSheet($SheetName, $HandoffFileName, $HandoffDatetime) {
    $ws = $wb.Worksheets.Item($SheetName)

    # Status: "Handoff transform failed" -> "Ready for handoff"
    $ws.Range("B2").Value = "Ready for handoff"

    # Latest Handoff File: new hyperlink cell pointing at the generated xlf
    # (styled to match the existing hyperlink cells in column A)
    $ws.Range("C2").Value = $HandoffFileName
    $ws.Hyperlinks.Add($ws.Range("C2"), "$repoBase/$HandoffFileName", "", "", $HandoffFileName)
    $ws.Range("C2").Font.Color = 15570276
    $ws.Range("C2").Font.Underline = 2
    $ws.Range("C2").Font.Name = "Calibri"
    $ws.Range("C2").Font.Size = 11

    # Latest Handoff Datetime
    $ws.Range("D2").Value = $HandoffDatetime

    # Handoff Reason: "Ignored" -> "Include"
    $ws.Range("H2").Value = "Include"
}

Update-LocaleSheet "zh-cn" "0e1993c0-23dc-4c91-9f34-364b5fb4460d.b76cd49c4b4099e031447afa5290f0430fe9375c.zh-cn.xlf" "2016-01-18 04:06:15"
Update-LocaleSheet "de-de" "0e1993c0-23dc-4c91-9f34-364b5fb4460d.b76cd49c4b4099e031447afa5290f0430fe9375c.de-de.xlf" "2016-01-18 04:06:30"
